# Daily attendance processing - 2026-01-24 15:01:28
# Reverse the order of the comma-separated "Recorded By" names in column G
# for every row whose list includes the "System" entry (alongside at least
# one other contributor). Rows with a single recorded-by value, or with
# multiple values that do NOT include "System", are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $parts = $text -split ", "

    if ($parts.Count -le 1) {
        continue
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p -eq "System") {
            $hasSystem = $true
        }
    }

    if (-not $hasSystem) {
        continue
    }

    $reversed = $parts[($parts.Count - 1)..0]
    $newText = $reversed -join ", "

    $cell.Value = $newText
}
